$wb = $excel.ActiveWorkbook

# --- Debts sheet: add an "active" column at the front ---
$wsDebts = $wb.Worksheets.Item("Debts")
[void]$wsDebts.Columns.Item(1).Insert()
$wsDebts.Cells.Item(1, 1).Value = "active"
[void]$wsDebts.Range("A1:A1048576").Select()

# --- Fixed Assets sheet: add an "active" column at the front ---
$wsFixedAssets = $wb.Worksheets.Item("Fixed Assets")
[void]$wsFixedAssets.Columns.Item(1).Insert()
$wsFixedAssets.Cells.Item(1, 1).Value = "active"

# Fixed Assets becomes the active tab, with C10 selected
[void]$wsFixedAssets.Activate()
[void]$wsFixedAssets.Range("C10").Select()
